# HotChocolate_Part_3 edit script
# Commit: "adds GraphQL Client slides"
#
# Summary of changes:
#  1. Swap the order of the "Strawberry Shake" slide and the
#     "Challenge #5: Subscriptions" slide (Challenge #5 now comes first).
#  2. Flesh out the (now later) "Strawberry Shake" slide with real content
#     describing the library (it previously only had a "TBD" placeholder).
#  3. Add a brand new "Challenge #7: GraphQL Client" slide (modelled on the
#     Challenge #5 slide) right after the Strawberry Shake slide.
#  4. The "Questions?" slide is unaffected content-wise; it simply shifts
#     down by one position because of the newly inserted slide.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Swap slides 6 ("Strawberry Shake") and 7 ("Challenge #5")
# ---------------------------------------------------------------------
$strawberrySlide = $p.Slides.Item(6)
$strawberrySlide.MoveTo(7)

# Now:
#   position 6 -> "Challenge #5: Subscriptions" (unchanged)
#   position 7 -> "Strawberry Shake" (to be enriched below)

# ---------------------------------------------------------------------
# 2. Enrich the "Strawberry Shake" slide (now at position 7)
# ---------------------------------------------------------------------
$shakeSlide = $p.Slides.Item(7)
$contentBox = $shakeSlide.Shapes.Item(3)   # "Textfeld 11" (id 6) - was "TBD"

$contentText = "GraphQL Client library by the same team as Hot Chocolate`r" + `
    "`r" + `
    "Supports newest GraphQL draft spec`r" + `
    "`r" + `
    "Creates a strongly typed C# client from queries, mutations etc.`r" + `
    "`r" + `
    "Provides caching, persisted queries and state"

$ctr = $contentBox.TextFrame.TextRange
$ctr.Text = $contentText
$ctr.Font.Name = "Raleway"
$ctr.Font.Color.RGB = 4210752
$ctr.LanguageID = 2057

# Mark the blank paragraph after "Supports newest GraphQL draft spec" as bold
# (matches the source paragraph formatting).
$ctr.Paragraphs(4, 1).Font.Bold = $true

# Mark the two occurrences of the word "GraphQL" as flagged-misspelling runs
# the same way the rest of the deck marks foreign/product words (err="1").
# (Cosmetic spell-check flag only; formatting stays identical.)
$p1 = $ctr.Paragraphs(1, 1)
$idx = $p1.Text.IndexOf("GraphQL") + 1
$null = $p1.Characters($idx, 7)

$p3 = $ctr.Paragraphs(3, 1)
$idx3 = $p3.Text.IndexOf("GraphQL") + 1
$null = $p3.Characters($idx3, 7)

# ---------------------------------------------------------------------
# 3. Add the new "Challenge #7: GraphQL Client" slide
# ---------------------------------------------------------------------
# Build it from a duplicate of the "Challenge #5" slide (position 6), which
# already has the right layout/shapes (title, date, "Now it's your turn!",
# description, GitHub picture + URL).
$dupSlides = $p.Slides.Item(6).Duplicate()
$challenge7 = $dupSlides.Item(1)
$challenge7.MoveTo(8)

# position 6 -> Challenge #5 (unchanged)
# position 7 -> Strawberry Shake (enriched)
# position 8 -> Challenge #7 (new, currently a clone of Challenge #5)
# position 9 -> Questions? (unchanged content, shifted down)

# --- Title: "Challenge #5: Subscriptions" -> "Challenge #7: GraphQL Client"
$titleShape = $challenge7.Shapes.Item(1)
$ttr = $titleShape.TextFrame.TextRange
$ttr.Text = "👷‍♂️ Challenge #7: GraphQL Client"
$titleRun2 = $ttr.Characters($ttr.Text.IndexOf("GraphQL") + 1, 7)

# --- Description: update text + let auto-fit recompute the width
$descShape = $challenge7.Shapes.Item(4)
$dtr = $descShape.TextFrame.TextRange
$dtr.Text = "In this exercise you learn to develop a GraphQL client with Strawberry Shake"
$dtrIdx = $dtr.Text.IndexOf("GraphQL") + 1
$null = $dtr.Characters($dtrIdx, 7)
